# Generate Report for Handoff
#
# - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   (Overview!E2/F2, zh-cn!C2, de-de!C2)
# - Timestamps bumped forward (report regenerated ~45-54s later):
#     Overview!G2 (Latest HO Xliff Generate Date)   2016-08-19 21:02:11 -> 2016-08-19 21:02:56
#     zh-cn!H2     (Latest Handoff Datetime)         2016-08-19 21:02:02 -> 2016-08-19 21:02:52
#     de-de!H2     (Latest Handoff Datetime)         2016-08-19 21:02:11 -> 2016-08-19 21:02:56
# - Status/date columns narrowed (Overview E & F, zh-cn C, de-de C)
#   from ~29.98 chars down to ~17.2 chars now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Value updates ----------------------------------------------------------

# Overview sheet: zh-cn / de-de status columns (E2, F2) and generate date (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 21:02:56"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 21:02:52"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 21:02:56"

# --- Column width updates ----------------------------------------------------
# Target on-disk width is ~17.2159881591797 characters. The Excel column-width
# model only persists widths in whole-pixel increments (characters = pixels/6
# for the default Calibri 11 font), so the closest reachable width is
# 17.16666... (98 px); feed a ColumnWidth comfortably inside that pixel's
# rounding bucket.
$narrowWidth = 16.3333333333333

$wsOverview.Range("E1").ColumnWidth = $narrowWidth
$wsOverview.Range("F1").ColumnWidth = $narrowWidth
$wsZhCn.Range("C1").ColumnWidth = $narrowWidth
$wsDeDe.Range("C1").ColumnWidth = $narrowWidth
